# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 79
$ws1.Range("F3").Value = 138
$ws1.Range("F5").Value = 374
$ws1.Range("F6").Value = 657
$ws1.Range("F7").Value = 107
$ws1.Range("F9").Value = 10822
$ws1.Range("F12").Value = 292
$ws1.Range("F13").Value = 206
$ws1.Range("F14").Value = 424
$ws1.Range("F15").Value = 9029
$ws1.Range("F17").Value = 734
$ws1.Range("F18").Value = 5307
$ws1.Range("F19").Value = 74
$ws1.Range("F20").Value = 3362
$ws1.Range("F21").Value = 2

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 79
$ws4.Range("F3").Value = 138
$ws4.Range("F5").Value = 374
$ws4.Range("F6").Value = 657
$ws4.Range("F8").Value = 107
$ws4.Range("F12").Value = 10822
$ws4.Range("F15").Value = 292
$ws4.Range("F16").Value = 206
$ws4.Range("F17").Value = 424
$ws4.Range("F18").Value = 9029
$ws4.Range("F20").Value = 734
$ws4.Range("F21").Value = 5307
$ws4.Range("F22").Value = 74
$ws4.Range("F23").Value = 3362
$ws4.Range("F24").Value = 2

$wb.Save()
